$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.756.91"
$ws.Range("E2").Value = "  +1.31%  "
$ws.Range("D3").Value = "2.280.79"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'119.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.27%  "
$ws.Range("D6").Value = "'267.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("D7").Value = "'0.650"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.19%  "
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Value = "'0.628"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.83%  "
$ws.Range("D10").Value = "'48.43"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.26%  "
$ws.Range("D11").Value = "'0.0948"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.78%  "
$ws.Range("D12").Value = "'9.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.98%  "
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").Value = "'15.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").Value = "'0.916"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.36%  "
$ws.Range("D16").Value = "2.624.84"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "2.278.40"
$ws.Range("E17").Value = "  +0.73%  "
$ws.Range("D18").Value = "43.744.65"
$ws.Range("E18").Value = "  +1.64%  "
$ws.Range("E19").Value = "  +3.48%  "
$ws.Range("D20").Value = "'6.95"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "'72.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.08%  "
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").Value = "'236.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("D24").Value = "'9.51"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.80%  "
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("D26").Value = "'12.03"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +6.72%  "
$ws.Range("E27").Value = "  +1.75%  "
$ws.Range("D28").Value = "'42.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.73%  "
$ws.Range("E29").Value = "  +0.81%  "
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("D31").Value = "'173.32"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.16%  "
$ws.Range("D32").Value = "'21.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("D33").Value = "'0.0926"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("D34").Value = "'5.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.12%  "
$ws.Range("D35").Value = "'0.131"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.15%  "
$ws.Range("E36").Value = "  +14.02%  "
$ws.Range("D37").Value = "'0.0388"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +11.39%  "
$ws.Range("E38").Value = "  +0.86%  "
$ws.Range("E39").Value = "  +5.43%  "
$ws.Range("D40").Value = "'2.56"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.95%  "
$ws.Range("B41").Value = "MultiversX"
$ws.Range("C41").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D41").Value = "'74.14"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").Value = "'13.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("E43").Value = "  +2.94%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("E45").Value = "  +2.69%  "
$ws.Range("D46").Value = "'5.92"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.97%  "
$ws.Range("D47").Value = "'74.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +42.49%  "
$ws.Range("E48").Value = "  +4.34%  "
$ws.Range("B49").Value = "Aave"
$ws.Range("C49").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D49").Value = "'102.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.62%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.101"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.94%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'8.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.27%  "
